$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values, keeping them stored as text (matching original inline-string formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.263.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.749.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.749.67'
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.28'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000250'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.374.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.763.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.251.48'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.730'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.893.71'
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.680.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.328'
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '424.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.797.87'
$ws.Range("D50").Style = "Normal"

# Update "Volume(1h)" column (E) values (percent text, already non-numeric so stays text automatically)
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  +0.56%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E21").Value = '  +11.47%  '
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("E24").Value = '  +5.60%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +3.17%  '
$ws.Range("E39").Value = '  +6.15%  '
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("E42").Value = '  +6.91%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("E51").Value = '  +0.58%  '
